# Adds author-search (UserName) and sort-validation (SortData) test-data sheets,
# and updates view selections on the existing sheets.

$wb = $excel.ActiveWorkbook

$jsonData = $wb.Worksheets.Item("JSONData")
$headers  = $wb.Worksheets.Item("Headers")

# --- Headers sheet: selection moves from A1:A4 to B1 -----------------------
[void]$headers.Range("B1").Select()

# --- JSONData sheet: no longer the tab-selected sheet; selection -> F1 -----
[void]$jsonData.Range("F1").Select()

# --- New sheet: UserName ----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$userName = $wb.Worksheets.Add($null, $lastSheet)
$userName.Name = "UserName"

# Cells that hold numeric-looking text must be forced to text so they are
# written as shared strings (matching "StatusCode"/"total_count" data),
# not as numbers.
$userName.Range("B2:B4").NumberFormat = "@"
$userName.Range("D2").NumberFormat = "@"
$userName.Range("D4").NumberFormat = "@"

$userName.Range("A1").Value = "URL"
$userName.Range("B1").Value = "StatusCode"
$userName.Range("C1").Value = "Status"
$userName.Range("D1").Value = "total_count"
$userName.Range("E1").Value = "Message"
$userName.Range("F1").Value = "/items[0]/owner/login"
$userName.Range("G1").Value = "/items[0]/owner/type"

$userName.Range("A2").Value = "?q=user:abishekk84"
$userName.Range("B2").Value = "200"
$userName.Range("C2").Value = "200 OK"
$userName.Range("D2").Value = "1"
$userName.Range("F2").Value = "Abishekk84"
$userName.Range("G2").Value = "User"

$userName.Range("A3").Value = "?q=user:123abc"
$userName.Range("B3").Value = "422"
$userName.Range("C3").Value = "422 Unprocessable Entity"
$userName.Range("E3").Value = "Validation Failed"

$userName.Range("A4").Value = "?q=user:123"
$userName.Range("B4").Value = "200"
$userName.Range("C4").Value = "200 OK"
$userName.Range("D4").Value = "0"

# formatting: reuse the same styles used on other sheets' header/data rows
$jsonData.Range("A1").Copy([System.Reflection.Missing]::Value)
$userName.Range("A1:G1").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

$jsonData.Range("A2").Copy([System.Reflection.Missing]::Value)
$userName.Range("A2:E4").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

$jsonData.Range("B2").Copy([System.Reflection.Missing]::Value)
$userName.Range("F2:G4").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

[void]$userName.Range("A1:A4").Select()

# --- New sheet: SortData -----------------------------------------------------
$sortData = $wb.Worksheets.Add($null, $userName)
$sortData.Name = "SortData"

$sortData.Range("A1").Value = "URL"
$sortData.Range("B1").Value = "SortParam"
$sortData.Range("C1").Value = "SortOrder"

$sortData.Range("A2").Value = "?q=api"

$sortData.Range("A3").Value = "?q=api"
$sortData.Range("B3").Value = "stars"
$sortData.Range("C3").Value = "desc"

$sortData.Range("A4").Value = "?q=api"
$sortData.Range("B4").Value = "stars"
$sortData.Range("C4").Value = "asc"

$jsonData.Range("A1").Copy([System.Reflection.Missing]::Value)
$sortData.Range("A1:C1").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

$jsonData.Range("A2").Copy([System.Reflection.Missing]::Value)
$sortData.Range("A2:A4").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

$jsonData.Range("B2").Copy([System.Reflection.Missing]::Value)
$sortData.Range("B2:C4").PasteSpecial(-4122, [System.Reflection.Missing]::Value)

$sortData.Activate()
